$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '58.803.35'
$ws.Range("E2").Value = '  +0.82%  '
$ws.Range("D3").Value = '2.625.88'
$ws.Range("E3").Value = '  +2.60%  '
$ws.Range("D5").Value = '520.54'
$ws.Range("E5").Value = '  +2.58%  '
$ws.Range("D6").Value = '145.17'
$ws.Range("E6").Value = '  -0.07%  '
$ws.Range("D7").Value = '0.998'
$ws.Range("E7").Value = '  -0.05%  '
$ws.Range("E8").Value = '  -0.70%  '
$ws.Range("D9").Value = '2.635.91'
$ws.Range("E9").Value = '  +2.54%  '
$ws.Range("D10").Value = '6.29'
$ws.Range("E10").Value = '  +0.34%  '
$ws.Range("E11").Value = '  +1.17%  '
$ws.Range("D12").Value = '0.334'
$ws.Range("E12").Value = '  +0.11%  '
$ws.Range("E13").Value = '  -0.90%  '
$ws.Range("D14").Value = '3.085.96'
$ws.Range("E14").Value = '  +2.55%  '
$ws.Range("D15").Value = '58.797.87'
$ws.Range("E15").Value = '  +0.83%  '
$ws.Range("D16").Value = '20.84'
$ws.Range("E16").Value = '  -1.00%  '
$ws.Range("E17").Value = '  -0.33%  '
$ws.Range("D18").Value = '2.633.83'
$ws.Range("E18").Value = '  +2.74%  '
$ws.Range("D19").Value = '345.52'
$ws.Range("E19").Value = '  +0.52%  '
$ws.Range("D20").Value = '4.47'
$ws.Range("E20").Value = '  -1.10%  '
$ws.Range("E21").Value = '  -0.53%  '
$ws.Range("D22").Value = '6.16'
$ws.Range("E22").Value = '  +2.09%  '
$ws.Range("E23").Value = '  -0.24%  '
$ws.Range("D24").Value = '61.49'
$ws.Range("E24").Value = '  +1.51%  '
$ws.Range("D25").Value = '0.416'
$ws.Range("E25").Value = '  -0.04%  '
$ws.Range("D26").Value = '0.164'
$ws.Range("E26").Value = '  +2.83%  '
$ws.Range("D27").Value = '0.996'
$ws.Range("E27").Value = '  -0.11%  '
$ws.Range("E28").Value = '  -1.71%  '
$ws.Range("D29").Value = '7.08'
$ws.Range("E29").Value = '  +1.60%  '
$ws.Range("E30").Value = '  -0.08%  '
$ws.Range("E31").Value = '  +3.50%  '
$ws.Range("D32").Value = '18.87'
$ws.Range("E32").Value = '  +0.91%  '
$ws.Range("E33").Value = '  +2.40%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '150.40'
$ws.Range("E34").Value = '  +0.53%  '
$ws.Range("D35").Value = '0.974'
$ws.Range("E35").Value = '  +3.65%  '
$ws.Range("D36").Value = '3.97'
$ws.Range("E36").Value = '  +0.38%  '
$ws.Range("E37").Value = '  +1.10%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '36.60'
$ws.Range("E38").Value = '  +1.58%  '
$ws.Range("D39").Value = '0.838'
$ws.Range("E39").Value = '  -1.25%  '
$ws.Range("E40").Value = '  +2.45%  '
$ws.Range("E41").Value = '  +1.41%  '
$ws.Range("B42").Value = 'FirstDigitalUSD'
$ws.Range("C42").Value = 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'
$ws.Range("D42").Value = '0.996'
$ws.Range("E42").Value = '  -0.02%  '
$ws.Range("B43").Value = 'Bittensor'
$ws.Range("C43").Value = 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'
$ws.Range("D43").Value = '276.12'
$ws.Range("E43").Value = '  -5.20%  '
$ws.Range("D44").Value = '0.0984'
$ws.Range("E44").Value = '  -0.94%  '
$ws.Range("D45").Value = '0.606'
$ws.Range("E45").Value = '  -0.03%  '
$ws.Range("D46").Value = '19.47'
$ws.Range("E46").Value = '  +2.25%  '
$ws.Range("E47").Value = '  -2.50%  '
$ws.Range("E48").Value = '  +0.43%  '
$ws.Range("D49").Value = '1.989.94'
$ws.Range("E49").Value = '  +3.71%  '
$ws.Range("E50").Value = '  +0.43%  '
$ws.Range("D51").Value = '4.63'
$ws.Range("E51").Value = '  +0.29%  '
